$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"236.7946826666667"
$ws.Range("H2").Value = [double]"710.384048"
$ws.Range("I2").Value = [double]"0.7123899543147419"
$ws.Range("J2").Value = [double]"0.7240508783182559"
$ws.Range("M2").Value = [double]"162.7225033333333"
$ws.Range("N2").Value = [double]"488.16751"
$ws.Range("O2").Value = [double]"0.5231437953541009"
$ws.Range("P2").Value = [double]"0.5247717033381212"
$ws.Range("Q2").Value = [double]"38531.82353954227"
$ws.Range("R2").Value = [double]"346786.4118558805"
$ws.Range("S2").Value = [double]"0.3726823844723487"
$ws.Range("T2").Value = [double]"0.3799614127185338"
$ws.Range("G3").Value = [double]"236.7946826666667"
$ws.Range("H3").Value = [double]"710.384048"
$ws.Range("I3").Value = [double]"0.7123899543147419"
$ws.Range("J3").Value = [double]"0.7240508783182559"
$ws.Range("O3").Value = [double]"0.0009322191998643353"
$ws.Range("P3").Value = [double]"0.0009351200601857102"
$ws.Range("Q3").Value = [double]"68.66201229631822"
$ws.Range("R3").Value = [double]"617.9581106668641"
$ws.Range("S3").Value = [double]"0.0006641035932026791"
$ws.Range("T3").Value = [double]"0.0006770745009104837"
$ws.Range("G4").Value = [double]"236.7946826666667"
$ws.Range("H4").Value = [double]"710.384048"
$ws.Range("I4").Value = [double]"0.7123899543147419"
$ws.Range("J4").Value = [double]"0.7240508783182559"
$ws.Range("M4").Value = [double]"61.580654"
$ws.Range("N4").Value = [double]"184.741962"
$ws.Range("O4").Value = [double]"0.1979783766474813"
$ws.Range("P4").Value = [double]"0.1985944416431287"
$ws.Range("Q4").Value = [double]"14581.9714223358"
$ws.Range("R4").Value = [double]"131237.7428010222"
$ws.Range("S4").Value = [double]"0.141037806695206"
$ws.Range("T4").Value = [double]"0.1437924799008309"
$ws.Range("G5").Value = [double]"236.7946826666667"
$ws.Range("H5").Value = [double]"710.384048"
$ws.Range("I5").Value = [double]"0.7123899543147419"
$ws.Range("J5").Value = [double]"0.7240508783182559"
$ws.Range("M5").Value = [double]"2.8947245"
$ws.Range("N5").Value = [double]"5.789449"
$ws.Range("O5").Value = [double]"0.009306378223129816"
$ws.Range("P5").Value = [double]"0.00622355841157717"
$ws.Range("Q5").Value = [double]"685.4553693849253"
$ws.Range("R5").Value = [double]"4112.732216309552"
$ws.Range("S5").Value = [double]"0.006629770357211158"
$ws.Range("T5").Value = [double]"0.00450617293416742"
$ws.Range("G6").Value = [double]"236.7946826666667"
$ws.Range("H6").Value = [double]"710.384048"
$ws.Range("I6").Value = [double]"0.7123899543147419"
$ws.Range("J6").Value = [double]"0.7240508783182559"
$ws.Range("M6").Value = [double]"83.559527"
$ws.Range("N6").Value = [double]"250.678581"
$ws.Range("O6").Value = [double]"0.2686392305754237"
$ws.Range("P6").Value = [double]"0.2694751765469873"
$ws.Range("Q6").Value = [double]"19786.45167974176"
$ws.Range("R6").Value = [double]"178078.0651176759"
$ws.Range("S6").Value = [double]"0.1913758891967735"
$ws.Range("T6").Value = [double]"0.1951137382638132"
$ws.Range("I7").Value = [double]"0.2358656137148928"
$ws.Range("J7").Value = [double]"0.2397264359793184"
$ws.Range("M7").Value = [double]"162.7225033333333"
$ws.Range("N7").Value = [double]"488.16751"
$ws.Range("O7").Value = [double]"0.5231437953541009"
$ws.Range("P7").Value = [double]"0.5247717033381212"
$ws.Range("Q7").Value = [double]"12757.52437504581"
$ws.Range("R7").Value = [double]"114817.7193754123"
$ws.Range("S7").Value = [double]"0.1233916323523333"
$ws.Range("T7").Value = [double]"0.125801650144044"
$ws.Range("I8").Value = [double]"0.2358656137148928"
$ws.Range("J8").Value = [double]"0.2397264359793184"
$ws.Range("O8").Value = [double]"0.0009322191998643353"
$ws.Range("P8").Value = [double]"0.0009351200601857102"
$ws.Range("S8").Value = [double]"0.0002198784536928077"
$ws.Range("T8").Value = [double]"0.000224172999241086"
$ws.Range("I9").Value = [double]"0.2358656137148928"
$ws.Range("J9").Value = [double]"0.2397264359793184"
$ws.Range("M9").Value = [double]"61.580654"
$ws.Range("N9").Value = [double]"184.741962"
$ws.Range("O9").Value = [double]"0.1979783766474813"
$ws.Range("P9").Value = [double]"0.1985944416431287"
$ws.Range("Q9").Value = [double]"4827.953591808653"
$ws.Range("R9").Value = [double]"43451.58232627787"
$ws.Range("S9").Value = [double]"0.04669629131023637"
$ws.Range("T9").Value = [double]"0.04760833770040997"
$ws.Range("I10").Value = [double]"0.2358656137148928"
$ws.Range("J10").Value = [double]"0.2397264359793184"
$ws.Range("M10").Value = [double]"2.8947245"
$ws.Range("N10").Value = [double]"5.789449"
$ws.Range("O10").Value = [double]"0.009306378223129816"
$ws.Range("P10").Value = [double]"0.00622355841157717"
$ws.Range("Q10").Value = [double]"226.9478259693622"
$ws.Range("R10").Value = [double]"1361.686955816173"
$ws.Range("S10").Value = [double]"0.002195054611061427"
$ws.Range("T10").Value = [double]"0.001491951477116503"
$ws.Range("I11").Value = [double]"0.2358656137148928"
$ws.Range("J11").Value = [double]"0.2397264359793184"
$ws.Range("M11").Value = [double]"83.559527"
$ws.Range("N11").Value = [double]"250.678581"
$ws.Range("O11").Value = [double]"0.2686392305754237"
$ws.Range("P11").Value = [double]"0.2694751765469873"
$ws.Range("Q11").Value = [double]"6551.10805594046"
$ws.Range("R11").Value = [double]"58959.97250346414"
$ws.Range("S11").Value = [double]"0.06336275698756888"
$ws.Range("T11").Value = [double]"0.06460032365850686"
$ws.Range("G12").Value = [double]"0.6305213333333334"
$ws.Range("H12").Value = [double]"1.891564"
$ws.Range("I12").Value = [double]"0.001896905195629352"
$ws.Range("J12").Value = [double]"0.001927955138422806"
$ws.Range("M12").Value = [double]"162.7225033333333"
$ws.Range("N12").Value = [double]"488.16751"
$ws.Range("O12").Value = [double]"0.5231437953541009"
$ws.Range("P12").Value = [double]"0.5247717033381212"
$ws.Range("Q12").Value = [double]"102.6000097650711"
$ws.Range("R12").Value = [double]"923.40008788564"
$ws.Range("S12").Value = [double]"0.0009923541834684524"
$ws.Range("T12").Value = [double]"0.001011736301949619"
$ws.Range("G13").Value = [double]"0.6305213333333334"
$ws.Range("H13").Value = [double]"1.891564"
$ws.Range("I13").Value = [double]"0.001896905195629352"
$ws.Range("J13").Value = [double]"0.001927955138422806"
$ws.Range("O13").Value = [double]"0.0009322191998643353"
$ws.Range("P13").Value = [double]"0.0009351200601857102"
$ws.Range("Q13").Value = [double]"0.1828286980724445"
$ws.Range("R13").Value = [double]"1.645458282652"
$ws.Range("S13").Value = [double]"1.768331443688095E-06"
$ws.Range("T13").Value = [double]"1.802869525077284E-06"
$ws.Range("G14").Value = [double]"0.6305213333333334"
$ws.Range("H14").Value = [double]"1.891564"
$ws.Range("I14").Value = [double]"0.001896905195629352"
$ws.Range("J14").Value = [double]"0.001927955138422806"
$ws.Range("M14").Value = [double]"61.580654"
$ws.Range("N14").Value = [double]"184.741962"
$ws.Range("O14").Value = [double]"0.1979783766474813"
$ws.Range("P14").Value = [double]"0.1985944416431287"
$ws.Range("Q14").Value = [double]"38.82791606761867"
$ws.Range("R14").Value = [double]"349.451244608568"
$ws.Range("S14").Value = [double]"0.000375546211284872"
$ws.Range("T14").Value = [double]"0.000382881174228078"
$ws.Range("G15").Value = [double]"0.6305213333333334"
$ws.Range("H15").Value = [double]"1.891564"
$ws.Range("I15").Value = [double]"0.001896905195629352"
$ws.Range("J15").Value = [double]"0.001927955138422806"
$ws.Range("M15").Value = [double]"2.8947245"
$ws.Range("N15").Value = [double]"5.789449"
$ws.Range("O15").Value = [double]"0.009306378223129816"
$ws.Range("P15").Value = [double]"0.00622355841157717"
$ws.Range("Q15").Value = [double]"1.825185551372667"
$ws.Range("R15").Value = [double]"10.951113308236"
$ws.Range("S15").Value = [double]"1.76533172039468E-05"
$ws.Range("T15").Value = [double]"1.199874141887468E-05"
$ws.Range("G16").Value = [double]"0.6305213333333334"
$ws.Range("H16").Value = [double]"1.891564"
$ws.Range("I16").Value = [double]"0.001896905195629352"
$ws.Range("J16").Value = [double]"0.001927955138422806"
$ws.Range("M16").Value = [double]"83.559527"
$ws.Range("N16").Value = [double]"250.678581"
$ws.Range("O16").Value = [double]"0.2686392305754237"
$ws.Range("P16").Value = [double]"0.2694751765469873"
$ws.Range("Q16").Value = [double]"52.68606437674267"
$ws.Range("R16").Value = [double]"474.174579390684"
$ws.Range("S16").Value = [double]"0.0005095831522283925"
$ws.Range("T16").Value = [double]"0.0005195360513011569"
$ws.Range("G17").Value = [double]"16.059769"
$ws.Range("H17").Value = [double]"32.119538"
$ws.Range("I17").Value = [double]"0.04831535056182164"
$ws.Range("J17").Value = [double]"0.032737474561192"
$ws.Range("M17").Value = [double]"162.7225033333333"
$ws.Range("N17").Value = [double]"488.16751"
$ws.Range("O17").Value = [double]"0.5231437953541009"
$ws.Range("P17").Value = [double]"0.5247717033381212"
$ws.Range("Q17").Value = [double]"2613.285814635063"
$ws.Range("R17").Value = [double]"15679.71488781038"
$ws.Range("S17").Value = [double]"0.02527587586677527"
$ws.Range("T17").Value = [double]"0.01717970028846514"
$ws.Range("G18").Value = [double]"16.059769"
$ws.Range("H18").Value = [double]"32.119538"
$ws.Range("I18").Value = [double]"0.04831535056182164"
$ws.Range("J18").Value = [double]"0.032737474561192"
$ws.Range("O18").Value = [double]"0.0009322191998643353"
$ws.Range("P18").Value = [double]"0.0009351200601857102"
$ws.Range("Q18").Value = [double]"4.656760211572333"
$ws.Range("R18").Value = [double]"27.940561269434"
$ws.Range("S18").Value = [double]"4.504049744190623E-05"
$ws.Range("T18").Value = [double]"3.061346918199002E-05"
$ws.Range("G19").Value = [double]"16.059769"
$ws.Range("H19").Value = [double]"32.119538"
$ws.Range("I19").Value = [double]"0.04831535056182164"
$ws.Range("J19").Value = [double]"0.032737474561192"
$ws.Range("M19").Value = [double]"61.580654"
$ws.Range("N19").Value = [double]"184.741962"
$ws.Range("O19").Value = [double]"0.1979783766474813"
$ws.Range("P19").Value = [double]"0.1985944416431287"
$ws.Range("Q19").Value = [double]"988.971078108926"
$ws.Range("R19").Value = [double]"5933.826468653556"
$ws.Range("S19").Value = [double]"0.009565394671383421"
$ws.Range("T19").Value = [double]"0.006501480481286055"
$ws.Range("G20").Value = [double]"16.059769"
$ws.Range("H20").Value = [double]"32.119538"
$ws.Range("I20").Value = [double]"0.04831535056182164"
$ws.Range("J20").Value = [double]"0.032737474561192"
$ws.Range("M20").Value = [double]"2.8947245"
$ws.Range("N20").Value = [double]"5.789449"
$ws.Range("O20").Value = [double]"0.009306378223129816"
$ws.Range("P20").Value = [double]"0.00622355841157717"
$ws.Range("Q20").Value = [double]"46.4886067886405"
$ws.Range("R20").Value = [double]"185.954427154562"
$ws.Range("S20").Value = [double]"0.0004496409263114198"
$ws.Range("T20").Value = [double]"0.0002037435851791001"
$ws.Range("G21").Value = [double]"16.059769"
$ws.Range("H21").Value = [double]"32.119538"
$ws.Range("I21").Value = [double]"0.04831535056182164"
$ws.Range("J21").Value = [double]"0.032737474561192"
$ws.Range("M21").Value = [double]"83.559527"
$ws.Range("N21").Value = [double]"250.678581"
$ws.Range("O21").Value = [double]"0.2686392305754237"
$ws.Range("P21").Value = [double]"0.2694751765469873"
$ws.Range("Q21").Value = [double]"1341.946701369263"
$ws.Range("R21").Value = [double]"8051.680208215578"
$ws.Range("S21").Value = [double]"0.01297939859990963"
$ws.Range("T21").Value = [double]"0.008821936737079719"
$ws.Range("G22").Value = [double]"0.5092873333333333"
$ws.Range("H22").Value = [double]"1.527862"
$ws.Range("I22").Value = [double]"0.001532176212914103"
$ws.Range("J22").Value = [double]"0.001557256002810873"
$ws.Range("M22").Value = [double]"162.7225033333333"
$ws.Range("N22").Value = [double]"488.16751"
$ws.Range("O22").Value = [double]"0.5231437953541009"
$ws.Range("P22").Value = [double]"0.5247717033381212"
$ws.Range("Q22").Value = [double]"82.87250979595777"
$ws.Range("R22").Value = [double]"745.85258816362"
$ws.Range("S22").Value = [double]"0.0008015484791751568"
$ws.Range("T22").Value = [double]"0.0008172038851285756"
$ws.Range("G23").Value = [double]"0.5092873333333333"
$ws.Range("H23").Value = [double]"1.527862"
$ws.Range("I23").Value = [double]"0.001532176212914103"
$ws.Range("J23").Value = [double]"0.001557256002810873"
$ws.Range("O23").Value = [double]"0.0009322191998643353"
$ws.Range("P23").Value = [double]"0.0009351200601857102"
$ws.Range("Q23").Value = [double]"0.1476751620851111"
$ws.Range("R23").Value = [double]"1.329076458766"
$ws.Range("S23").Value = [double]"1.428324083253953E-06"
$ws.Range("T23").Value = [double]"1.456221327073062E-06"
$ws.Range("G24").Value = [double]"0.5092873333333333"
$ws.Range("H24").Value = [double]"1.527862"
$ws.Range("I24").Value = [double]"0.001532176212914103"
$ws.Range("J24").Value = [double]"0.001557256002810873"
$ws.Range("M24").Value = [double]"61.580654"
$ws.Range("N24").Value = [double]"184.741962"
$ws.Range("O24").Value = [double]"0.1979783766474813"
$ws.Range("P24").Value = [double]"0.1985944416431287"
$ws.Range("Q24").Value = [double]"31.36224706058267"
$ws.Range("R24").Value = [double]"282.260223545244"
$ws.Range("S24").Value = [double]"0.0003033377593706198"
$ws.Range("T24").Value = [double]"0.0003092623863736357"
$ws.Range("G25").Value = [double]"0.5092873333333333"
$ws.Range("H25").Value = [double]"1.527862"
$ws.Range("I25").Value = [double]"0.001532176212914103"
$ws.Range("J25").Value = [double]"0.001557256002810873"
$ws.Range("M25").Value = [double]"2.8947245"
$ws.Range("N25").Value = [double]"5.789449"
$ws.Range("O25").Value = [double]"0.009306378223129816"
$ws.Range("P25").Value = [double]"0.00622355841157717"
$ws.Range("Q25").Value = [double]"1.474246521339667"
$ws.Range("R25").Value = [double]"8.845479128038001"
$ws.Range("S25").Value = [double]"1.425901134186132E-05"
$ws.Range("T25").Value = [double]"1.199874141887468E-05"
$ws.Range("G26").Value = [double]"0.5092873333333333"
$ws.Range("H26").Value = [double]"1.527862"
$ws.Range("I26").Value = [double]"0.001532176212914103"
$ws.Range("J26").Value = [double]"0.001557256002810873"
$ws.Range("M26").Value = [double]"83.559527"
$ws.Range("N26").Value = [double]"250.678581"
$ws.Range("O26").Value = [double]"0.2686392305754237"
$ws.Range("P26").Value = [double]"0.2694751765469873"
$ws.Range("Q26").Value = [double]"42.55580868042467"
$ws.Range("R26").Value = [double]"383.002278123822"
$ws.Range("S26").Value = [double]"0.0004116026389432111"
$ws.Range("T26").Value = [double]"0.0004196418362863156"
